$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New word-tracking columns (E/F/G) -- mark which source list each word is new to.
$cells = @(
  "E4",
  "E6","F6",
  "F10","G10",
  "E12","F12",
  "F13",
  "E14","F14",
  "E16","F16","G16",
  "E17","F17",
  "E20",
  "E25","G25",
  "E29","F29",
  "E33","F33",
  "F34",
  "E37","F37",
  "E38","F38","G38"
)
foreach ($c in $cells) {
  $ws.Range($c).Value = 1
}

# Column B/C (the Spanish translation helper columns) are hidden now; column D widened.
$ws.Columns("B").Hidden = $true
$ws.Columns("C").Hidden = $true
$ws.Columns("D").ColumnWidth = 24.022135416666668

# Update the view: selection moves to G23, scroll position to A10 while keeping
# the header row (row 1) frozen.
$ws.Range("G23").Select()
